$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.902.61"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.442.93"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.72"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.25"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +2.08%  "

$ws.Range("E9").Value = "  +11.32%  "

$ws.Range("E10").Value = "  -1.63%  "

$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("E12").Value = "  -5.05%  "

$ws.Range("E13").Value = "  +4.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.794.16"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.886.58"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.19"
$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.437.95"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.52"
$ws.Range("E18").Value = "  -1.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.93"
$ws.Range("E19").Value = "  +1.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.95"
$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("E21").Value = "  +1.51%  "

$ws.Range("E22").Value = "  +2.80%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.03"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.70"
$ws.Range("E25").Value = "  +2.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.558.28"
$ws.Range("E26").Value = "  -1.13%  "

$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("E28").Value = "  +1.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0819"
$ws.Range("E29").Value = "  +0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.14"
$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "428.94"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("E33").Value = "  +2.19%  "

$ws.Range("E34").Value = "  +0.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.45"
$ws.Range("E35").Value = "  +0.37%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("E38").Value = "  +1.50%  "

$ws.Range("E39").Value = "  -2.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.298"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("E41").Value = "  +3.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.35"
$ws.Range("E42").Value = "  -1.18%  "

$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.05"
$ws.Range("E44").Value = "  -0.62%  "

$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.70"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0716"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.482"
$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0924"
$ws.Range("E50").Value = "  +1.77%  "

$ws.Range("E51").Value = "  +1.13%  "
